$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameter_values")

# Insert two new rows to make room for the new postpartum-haemorrhage parameters
$ws.Rows.Item(34).Insert()
$ws.Rows.Item(48).Insert()

$ws.Range("A34").Value = 'severity_maternal_haemorrhage'
$ws.Range("B34").Value = '[0.3, 0.7]'
$ws.Range("C34").Value = 'Dummy'
$ws.Range("A35").Value = 'prob_cord_prolapse'
$ws.Range("B35").Value = 0.004
$ws.Range("A36").Value = 'cfr_aph'
$ws.Range("B36").Value = 0.6
$ws.Range("C36").Value = 0.02
$ws.Range("A37").Value = 'cfr_eclampsia'
$ws.Range("B37").Value = 0.5
$ws.Range("C37").Value = 0.184
$ws.Range("A38").Value = 'cfr_sepsis'
$ws.Range("B38").Value = 0.5
$ws.Range("C38").Value = 0.33
$ws.Range("A39").Value = 'cfr_uterine_rupture'
$ws.Range("B39").Value = 0.8
$ws.Range("C39").Value = 0.345
$ws.Range("A40").Value = 'prob_still_birth_obstructed_labour'
$ws.Range("B40").Value = 0.38
$ws.Range("A41").Value = 'prob_still_birth_antepartum_haem'
$ws.Range("B41").Value = 0.38
$ws.Range("A42").Value = 'prob_still_birth_sepsis'
$ws.Range("B42").Value = 0.25
$ws.Range("A43").Value = 'prob_still_birth_uterine_rupture'
$ws.Range("B43").Value = 0.93
$ws.Range("A44").Value = 'prob_still_birth_eclampsia'
$ws.Range("B44").Value = 0.03
$ws.Range("A45").Value = 'prob_pp_eclampsia'
$ws.Range("B45").Value = 0.01
$ws.Range("A46").Value = 'prob_pph'
$ws.Range("B46").Value = 0.03
$ws.Range("A47").Value = 'rr_pph_pl_ol'
$ws.Range("B47").Value = 5
$ws.Range("A48").Value = 'prob_pph_source'
$ws.Range("B48").Value = '[0.67, 0.33]'
$ws.Range("A49").Value = 'prob_pp_sepsis'
$ws.Range("B49").Value = 0.05
$ws.Range("A50").Value = 'cfr_pp_pph'
$ws.Range("B50").Value = 0.5
$ws.Range("C50").Value = 0.1
$ws.Range("A51").Value = 'cfr_pp_eclampsia'
$ws.Range("B51").Value = 0.5
$ws.Range("C51").Value = 0.184
$ws.Range("A52").Value = 'cfr_pp_sepsis'
$ws.Range("B52").Value = 0.5
$ws.Range("C52").Value = 0.33
$ws.Range("A53").Value = 'prob_neonatal_sepsis'
$ws.Range("B53").Value = 0.15
$ws.Range("A54").Value = 'prob_neonatal_birth_asphyxia'
$ws.Range("B54").Value = 0.16
$ws.Range("A55").Value = 'odds_homebirth'
$ws.Range("B55").Value = 0.5
$ws.Range("A56").Value = 'or_homebirth_unmarried'
$ws.Range("B56").Value = 1.83
$ws.Range("A57").Value = 'or_homebirth_wealth_4'
$ws.Range("B57").Value = 0.51
$ws.Range("A58").Value = 'or_homebirth_wealth_5'
$ws.Range("B58").Value = 0.43
$ws.Range("A59").Value = 'or_homebirth_urban'
$ws.Range("B59").Value = 0.39
$ws.Range("A60").Value = 'prob_successful_induction'
$ws.Range("B60").Value = 0.761
$ws.Range("A61").Value = 'rr_maternal_sepsis_clean_delivery'
$ws.Range("B61").Value = 0.7
$ws.Range("A62").Value = 'rr_newborn_sepsis_clean_delivery'
$ws.Range("B62").Value = 0.7
$ws.Range("A63").Value = 'rr_sepsis_post_abx_prom'
$ws.Range("B63").Value = 0.7
$ws.Range("A64").Value = 'rr_sepsis_post_abx_pprom'
$ws.Range("B64").Value = 0.7
$ws.Range("A65").Value = 'rr_newborn_sepsis_proph_abx'
$ws.Range("B65").Value = 0.8
$ws.Range("A66").Value = 'rr_pph_amtsl'
$ws.Range("B66").Value = 0.34
$ws.Range("A67").Value = 'prob_cure_antibiotics'
$ws.Range("B67").Value = 0.5
$ws.Range("A68").Value = 'prob_cure_mgso4'
$ws.Range("B68").Value = 0.57
$ws.Range("A69").Value = 'prob_prevent_mgso4'
$ws.Range("B69").Value = 0.41
$ws.Range("A70").Value = 'prob_cure_diazepam'
$ws.Range("B70").Value = 0.8
$ws.Range("A71").Value = 'prob_cure_blood_transfusion'
$ws.Range("B71").Value = 0.2
$ws.Range("A72").Value = 'prob_cure_oxytocin'
$ws.Range("B72").Value = 0.5
$ws.Range("A73").Value = 'prob_cure_misoprostol'
$ws.Range("B73").Value = 0.3
$ws.Range("A74").Value = 'prob_cure_uterine_massage'
$ws.Range("B74").Value = 0.15
$ws.Range("A75").Value = 'prob_cure_uterine_tamponade'
$ws.Range("B75").Value = 0.6
$ws.Range("A76").Value = 'prob_cure_uterine_ligation'
$ws.Range("B76").Value = 0.8
$ws.Range("A77").Value = 'prob_cure_b_lynch'
$ws.Range("B77").Value = 0.8
$ws.Range("A78").Value = 'prob_cure_hysterectomy'
$ws.Range("B78").Value = 0.95
$ws.Range("A79").Value = 'prob_successful_manual_removal_placenta'
$ws.Range("B79").Value = 0.75
$ws.Range("A80").Value = 'prob_cure_uterine_repair'
$ws.Range("B80").Value = 0.7
$ws.Range("A81").Value = 'prob_successful_assisted_vaginal_delivery'
$ws.Range("B81").Value = 0.7
$ws.Range("A82").Value = 'dummy_prob_health_centre'
$ws.Range("B82").Value = 0.7
$ws.Range("C82").Value = 'DUMMY'
$ws.Range("A83").Value = 'squeeze_factor_threshold_delivery_attendance'
$ws.Range("B83").Value = 0.8
$ws.Range("C83").Value = 'DUMMY'
$ws.Range("A84").Value = 'squeeze_factor_threshold_sba_did_not_run'
$ws.Range("B84").Value = 0.9
$ws.Range("C84").Value = 'DUMMY'
$ws.Range("A85").Value = 'sensitivity_of_assessment_of_obstructed_labour_hc'
$ws.Range("B85").Value = 0.25
$ws.Range("C85").Value = 'DUMMY'
$ws.Range("A86").Value = 'sensitivity_of_assessment_of_obstructed_labour_hp'
$ws.Range("B86").Value = 0.5
$ws.Range("C86").Value = 'DUMMY'
$ws.Range("A87").Value = 'sensitivity_of_assessment_of_obstructed_labour_for_cs'
$ws.Range("B87").Value = 0.8
$ws.Range("C87").Value = 'DUMMY'
$ws.Range("A88").Value = 'sensitivity_of_assessment_of_sepsis_hc'
$ws.Range("B88").Value = 0.5
$ws.Range("C88").Value = 'DUMMY'
$ws.Range("A89").Value = 'sensitivity_of_assessment_of_sepsis_hp'
$ws.Range("B89").Value = 0.25
$ws.Range("C89").Value = 'DUMMY'
$ws.Range("A90").Value = 'sensitivity_of_assessment_of_hypertension_hc'
$ws.Range("B90").Value = 0.4
$ws.Range("C90").Value = 'DUMMY'
$ws.Range("A91").Value = 'sensitivity_of_assessment_of_hypertension_hp'
$ws.Range("B91").Value = 0.8
$ws.Range("C91").Value = 'DUMMY'
$ws.Range("A92").Value = 'sensitivity_of_assessment_of_severe_pe_hc'
$ws.Range("B92").Value = 0.3
$ws.Range("C92").Value = 'DUMMY'
$ws.Range("A93").Value = 'sensitivity_of_assessment_of_severe_pe_hp'
$ws.Range("B93").Value = 0.6
$ws.Range("C93").Value = 'DUMMY'
$ws.Range("A94").Value = 'sensitivity_of_referral_assessment_of_antepartum_haem_hc'
$ws.Range("B94").Value = 0.4
$ws.Range("C94").Value = 'DUMMY'
$ws.Range("A95").Value = 'sensitivity_of_treatment_assessment_of_antepartum_haem_hp'
$ws.Range("B95").Value = 0.8
$ws.Range("C95").Value = 'DUMMY'
$ws.Range("A96").Value = 'sensitivity_of_referral_assessment_of_uterine_rupture_hc'
$ws.Range("B96").Value = 0.5
$ws.Range("C96").Value = 'DUMMY'
$ws.Range("A97").Value = 'sensitivity_of_treatment_assessment_of_uterine_rupture_hp'
$ws.Range("B97").Value = 0.8
$ws.Range("C97").Value = 'DUMMY'

# Update view state to match target (scroll position / selection)
$ws.Application.ActiveWindow.ScrollRow = 67
$ws.Range("A79").Select()

$ws2 = $wb.Worksheets.Item("parameter_sources")
$ws2.Application.ActiveWindow.ScrollRow = 4